$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to Text format so numeric-looking
# strings like "569.18" are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '66.116.70'
$ws.Range("E2").Value = '  -4.66%  '
$ws.Range("D3").Value = '3.298.13'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '569.18'
$ws.Range("E5").Value = '  -3.49%  '
$ws.Range("D6").Value = '179.97'
$ws.Range("E6").Value = '  -5.05%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.597'
$ws.Range("E8").Value = '  -1.46%  '
$ws.Range("D9").Value = '3.293.75'
$ws.Range("E9").Value = '  -0.34%  '
$ws.Range("D10").Value = '0.128'
$ws.Range("E10").Value = '  -3.65%  '
$ws.Range("D11").Value = '6.61'
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("D12").Value = '0.401'
$ws.Range("E12").Value = '  -3.84%  '
$ws.Range("D13").Value = '3.869.57'
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("E15").Value = '  -3.85%  '
$ws.Range("D16").Value = '66.201.70'
$ws.Range("E16").Value = '  -4.42%  '
$ws.Range("E17").Value = '  -3.46%  '
$ws.Range("D18").Value = '3.316.64'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = '433.92'
$ws.Range("E19").Value = '  +3.20%  '
$ws.Range("E20").Value = '  -0.64%  '
$ws.Range("D21").Value = '5.65'
$ws.Range("E21").Value = '  -2.59%  '
$ws.Range("D22").Value = '7.61'
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("D23").Value = '73.35'
$ws.Range("E23").Value = '  +1.81%  '
$ws.Range("D24").Value = '0.997'
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("D25").Value = '3.453.41'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").Value = '0.510'
$ws.Range("E26").Value = '  -1.05%  '
$ws.Range("D27").Value = '0.0000117'
$ws.Range("E27").Value = '  -2.51%  '
$ws.Range("D28").Value = '0.190'
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").Value = '8.92'
$ws.Range("E29").Value = '  -7.25%  '
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("E31").Value = '  -1.36%  '
$ws.Range("D32").Value = '22.70'
$ws.Range("E32").Value = '  -1.31%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("D34").Value = '5.25'
$ws.Range("E34").Value = '  -5.72%  '
$ws.Range("D35").Value = '6.71'
$ws.Range("E35").Value = '  -3.90%  '
$ws.Range("E36").Value = '  -4.87%  '
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("D38").Value = '159.43'
$ws.Range("E38").Value = '  -2.99%  '
$ws.Range("D39").Value = '27.13'
$ws.Range("E39").Value = '  +1.04%  '
$ws.Range("D40").Value = '1.83'
$ws.Range("E40").Value = '  -5.37%  '
$ws.Range("D41").Value = '2.773.24'
$ws.Range("E41").Value = '  +2.53%  '
$ws.Range("D42").Value = '0.781'
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("D43").Value = '4.43'
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("D44").Value = '6.12'
$ws.Range("E44").Value = '  -4.86%  '
$ws.Range("D45").Value = '0.0672'
$ws.Range("E45").Value = '  -2.43%  '
$ws.Range("D46").Value = '40.07'
$ws.Range("E46").Value = '  -1.79%  '
$ws.Range("D47").Value = '24.06'
$ws.Range("E47").Value = '  -4.37%  '
$ws.Range("D48").Value = '2.33'
$ws.Range("E48").Value = '  -6.03%  '
$ws.Range("D49").Value = '317.42'
$ws.Range("E49").Value = '  -7.01%  '
$ws.Range("E50").Value = '  -3.44%  '
$ws.Range("D51").Value = '0.972'
$ws.Range("E51").Value = '  -2.98%  '

# Restore the original (default/General) formatting so the workbook
# styling matches the source - only the text values should differ.
$priceRange.ClearFormats()

